# 29 Jan 2025 - LV Contacts - Final
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# Update the CompanyName values for the two contact rows.
$ws.Range("A2").Value = "Houlihan Lokey - NY"
$ws.Range("A3").Value = "Houlihan Lokey - NY"

# Re-fit column A now that the text is longer.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Move the selection/view back to the top-left, selecting B9 as in the saved file.
$ws.Activate()
$ws.Range("B9").Select()
